# market_health_data.xlsx refresh: 2025-11-05 09:41 AM
# A new "Stock List" snapshot row (CAPTRU-RE1) is prepended to the table;
# every existing data row shifts down by one, and the oldest row (which
# would fall past the fixed A1:H76 table range) drops off the bottom.

$wb = $excel.ActiveWorkbook

# --- 1. Bump the "Last Updated" timestamp on the Metadata sheet ---------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Cells.Item(2, 1).Value = "05 Nov 2025, 09:41 AM"

# --- 2. Push the "Stock List" table down by one row and insert the -------
#        newest entry (CAPTRU-RE1) at the top -----------------------------
$ws = $wb.Worksheets.Item("Stock List")

# Columns B:E hold Stock / Stock Name / Price / % Change - shift rows
# 2..75 down into rows 3..76 (row 76's original contents fall off the end).
$nameBlock = $ws.Range("B2:E75").Value()
$ws.Range("B3:E76").Value = $nameBlock

# Column H holds Market Cap - shift the same way (columns F/G are static
# "N/A" placeholders and are left untouched).
$capBlock = $ws.Range("H2:H75").Value()
$ws.Range("H3:H76").Value = $capBlock

# Write the new first row: CAPTRU-RE1
$ws.Range("B2:C2").Value = "CAPTRU-RE1"
$ws.Cells.Item(2, 4).Value = 5.67
$ws.Cells.Item(2, 5).Value = -11.9565
